# Update cryptos list with the latest coinranking.com snapshot data.
# Prices (D) and 1h volume deltas (E) are refreshed; a couple of rows
# (39/40 and 44/45/46) also swap rank position, changing Coin (B) and
# Link (C) along with their Price/Volume.
#
# Source cells are stored as plain text, and several "Price" values
# look like numbers (e.g. 1.003, 0.4484). A leading apostrophe forces
# Excel to keep them as text instead of auto-converting to a number,
# which would silently change digits such as trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.438.45"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.830.42"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'330.22"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4484"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("D8").Value = "'0.3796"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "'44.75"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").Value = "'0.07806"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").Value = "'1.145"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "'22.57"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "'6.389"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "'7.537"
$ws.Range("D16").Value = "1.839.72"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").Value = "'93.89"
$ws.Range("E17").Value = "  +15.91%  "
$ws.Range("D18").Value = "'0.00001089"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "'0.06398"
$ws.Range("E19").Value = "  -4.39%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'17.60"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'6.392"
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("D23").Value = "'0.5432"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "28.494.15"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "'11.76"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "'2.287"
$ws.Range("E26").Value = "  -6.14%  "
$ws.Range("D27").Value = "'20.92"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").Value = "'154.10"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").Value = "2.045.35"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'129.30"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'1.215"
$ws.Range("E32").Value = "  -6.44%  "
$ws.Range("D33").Value = "'5.924"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").Value = "'0.09330"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "'3.674"
$ws.Range("E35").Value = "  -7.43%  "
$ws.Range("D36").Value = "'13.01"
$ws.Range("E36").Value = "  +7.35%  "
$ws.Range("D37").Value = "'0.02369"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "'0.2209"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06326"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6685"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "'5.237"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'8.217"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("D43").Value = "'1.198"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.0000"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.407"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.07"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'0.6183"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "'3.784"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "'2.064"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").Value = "'127.75"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  -0.20%  "
